$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (row 1)
$ws.Range("A1").Value = "MIGRATION DATE"
$ws.Range("B1").Value = "FINANCIAL INSTITUTION NAME"
$ws.Range("C1").Value = "ENTITY ID"
$ws.Range("D1").Value = "ADDRESS"

# New data row (row 2)
# Force the migration date to be stored as literal text (not auto-converted
# to a date serial number) by temporarily marking the cell as Text before
# entering the value, then clearing that temporary formatting again so the
# cell is left with the workbook's default (unstyled) formatting.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-10-17"
$ws.Range("A2").ClearFormats()

$ws.Range("B2").Value = "ZZZ"
$ws.Range("C2").Value = "456CDX009"
$ws.Range("D2").Value = "Anna Nagar"

# The header cells A1/B1 already carried the bold + thin-border + centered
# style from the original workbook. Extend that same style to the two new
# header cells (C1/D1) by copying the format from A1, so the whole header
# row shares one consistent style (matches A1/B1 rather than creating a
# brand-new style entry).
$ws.Range("A1").Copy() | Out-Null
$ws.Range("C1:D1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
